$d = $word.ActiveDocument

# Word "BGR" integer equivalents of the OOXML hex colors used in the target:
#   00CC33 (green) -> 0x33CC00 -> 3394560
#   FF3333 (red)   -> 0x3333FF -> 3355647
$greenColor = 3394560
$redColor   = 3355647

# --- Paragraph 5: "Many figures showing data with error bars ..." -------
# Split into a green run (original sentence + " - ") and a new red run
# ("but maybe need to define better.").
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5.Font.Color = $greenColor
$r5.InsertAfter(" - ")

$p5b = $d.Paragraphs.Item(5)
$p5b.Range.InsertAfter("but maybe need to define better.")

$p5c = $d.Paragraphs.Item(5)
$findRange = $p5c.Range.Duplicate
$findRange.Find.Execute("but maybe need to define better.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRange.Font.Color = $redColor

# --- Remaining paragraphs: simply recolor the whole paragraph green -----
$greenParagraphs = @(30, 31, 32, 33, 34, 35, 36, 38, 41, 42, 43, 44, 46, 49)
foreach ($idx in $greenParagraphs) {
    $para = $d.Paragraphs.Item($idx)
    $para.Range.Font.Color = $greenColor
}
